# Generate Report for Handoff
# - Bumps the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" stamps
#   for the e2e test set from 2016-08-25 18:32:13 to 2016-08-25 18:32:35 (Overview
#   sheet + the de-de handoff sheet, which shared the same timestamp string).
# - Bumps the zh-cn "Latest Handoff Datetime" stamp from 18:31:58 to 18:32:30.
# - Flips the Priority column for rows 4-7 (the .md-sourced files) from "low" to "ht"
#   on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G), rows 4-7.
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-08-25 18:32:35"
}

# zh-cn sheet: Priority (E) and Latest Handoff Datetime (H), rows 4-7.
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-25 18:32:30"
}

# de-de sheet: Priority (E) and Latest Handoff Datetime (H), rows 4-7.
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-25 18:32:35"
}
